$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5963.4736
$ws.Range("I40").Value = 4556.222
$ws.Range("K40").Value = 4556.222
$ws.Range("M40").Value = -4381.222
$ws.Range("H88").Value = 901977.8
$ws.Range("I88").Value = 2251934.8
$ws.Range("K88").Value = 2251934.8
$ws.Range("M88").Value = -2251528.8
$ws.Range("H91").Value = 901977.8
$ws.Range("I91").Value = 2251934.8
$ws.Range("K91").Value = 2251934.8
$ws.Range("M91").Value = -2250530.8
$ws.Range("I92").Value = 7143090.5
$ws.Range("K92").Value = 7143090.5
$ws.Range("M92").Value = -7141842.5
$ws.Range("H137").Value = 728078.5600000001
$ws.Range("I137").Value = 2454.5334
$ws.Range("J137").Value = 2904950.5
$ws.Range("K137").Value = 7363.600199999999
$ws.Range("L137").Value = 8714851.5
$ws.Range("M137").Value = -4813.600199999999
$ws.Range("N137").Value = -8719951.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4123.381
$ws.Range("I32").Value = 1610.4412
$ws.Range("K32").Value = 1610.4412
$ws.Range("M32").Value = -1323.4412
$ws.Range("H63").Value = 1810.8334
$ws.Range("I63").Value = 1761
$ws.Range("K63").Value = 1761
$ws.Range("M63").Value = -1075
$ws.Range("H66").Value = 1810.8334
$ws.Range("I66").Value = 1761
$ws.Range("K66").Value = 8805
$ws.Range("M66").Value = -5373
$ws.Range("H102").Value = 113867.1
$ws.Range("I102").Value = 146487.14
$ws.Range("K102").Value = 146487.14
$ws.Range("M102").Value = -144865.14
$ws.Range("H110").Value = 1737.5
$ws.Range("I110").Value = 1700
$ws.Range("K110").Value = 1700
$ws.Range("M110").Value = 345
$ws.Range("H132").Value = 1082.9286
$ws.Range("I132").Value = 721.75
$ws.Range("K132").Value = 2165.25
$ws.Range("M132").Value = 364.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1340.6666
$ws.Range("I20").Value = 1385
$ws.Range("J20").Value = 1207.6666
$ws.Range("K20").Value = 1385
$ws.Range("L20").Value = 1207.6666
$ws.Range("M20").Value = -1138
$ws.Range("N20").Value = -1701.6666
$ws.Range("H86").Value = 3376.3333
$ws.Range("I86").Value = 2599.6
$ws.Range("K86").Value = 2599.6
$ws.Range("M86").Value = -1476.6
$ws.Range("H89").Value = 3376.3333
$ws.Range("I89").Value = 2599.6
$ws.Range("K89").Value = 12998
$ws.Range("M89").Value = -7382
$ws.Range("H94").Value = 2164.88
$ws.Range("I94").Value = 2315
$ws.Range("J94").Value = 1564.4
$ws.Range("K94").Value = 2315
$ws.Range("L94").Value = 1564.4
$ws.Range("M94").Value = -1864
$ws.Range("N94").Value = -2466.4
$ws.Range("H105").Value = 128874.625
$ws.Range("I105").Value = 252999.25
$ws.Range("K105").Value = 252999.25
$ws.Range("M105").Value = -251252.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 997.5
$ws.Range("J23").Value = 997.5
$ws.Range("L23").Value = 997.5
$ws.Range("N23").Value = -1477.5
$ws.Range("H27").Value = 997.5
$ws.Range("J27").Value = 997.5
$ws.Range("L27").Value = 997.5
$ws.Range("N27").Value = -1381.5
$ws.Range("H132").Value = 1979029.6
$ws.Range("I132").Value = 2676452.5
$ws.Range("K132").Value = 8029357.5
$ws.Range("M132").Value = -8026827.5
$ws.Range("H134").Value = 2589924.2
$ws.Range("I134").Value = 3762963
$ws.Range("K134").Value = 11288889
$ws.Range("M134").Value = -11286354

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 11028.4
$ws.Range("I18").Value = 13285.5
$ws.Range("K18").Value = 39856.5
$ws.Range("M18").Value = -39687.5
$ws.Range("H137").Value = 5840.65
$ws.Range("J137").Value = 10919
$ws.Range("L137").Value = 32757
$ws.Range("N137").Value = -42957

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 75321.30499999999
$ws.Range("I70").Value = 7929.7144
$ws.Range("J70").Value = 153944.83
$ws.Range("K70").Value = 7929.7144
$ws.Range("L70").Value = 153944.83
$ws.Range("M70").Value = -7659.7144
$ws.Range("N70").Value = -154484.83
$ws.Range("H73").Value = 75321.30499999999
$ws.Range("I73").Value = 7929.7144
$ws.Range("J73").Value = 153944.83
$ws.Range("K73").Value = 7929.7144
$ws.Range("L73").Value = 153944.83
$ws.Range("M73").Value = -6993.7144
$ws.Range("N73").Value = -155816.83
$ws.Range("H102").Value = 2939.2727
$ws.Range("I102").Value = 2939.2727
$ws.Range("K102").Value = 2939.2727
$ws.Range("M102").Value = -1317.2727
$ws.Range("H132").Value = 4933.5293
$ws.Range("I132").Value = 4052
$ws.Range("K132").Value = 12156
$ws.Range("M132").Value = -9626

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3069.64
$ws.Range("I7").Value = 2176.9285
$ws.Range("J7").Value = 4205.8184
$ws.Range("K7").Value = 2176.9285
$ws.Range("L7").Value = 4205.8184
$ws.Range("M7").Value = -2064.9285
$ws.Range("N7").Value = -4429.8184
$ws.Range("H22").Value = 15549.857
$ws.Range("I22").Value = 1750
$ws.Range("J22").Value = 25899.75
$ws.Range("K22").Value = 1750
$ws.Range("L22").Value = 25899.75
$ws.Range("M22").Value = -1455
$ws.Range("N22").Value = -26489.75
$ws.Range("H27").Value = 15549.857
$ws.Range("I27").Value = 1750
$ws.Range("J27").Value = 25899.75
$ws.Range("K27").Value = 1750
$ws.Range("L27").Value = 25899.75
$ws.Range("M27").Value = -1643
$ws.Range("N27").Value = -26113.75
$ws.Range("I40").Value = 3680.8
$ws.Range("K40").Value = 3680.8
$ws.Range("M40").Value = -3544.8
$ws.Range("H46").Value = 5926.48
$ws.Range("I46").Value = 9636.75
$ws.Range("J46").Value = 2501.6155
$ws.Range("K46").Value = 9636.75
$ws.Range("L46").Value = 2501.6155
$ws.Range("M46").Value = -9448.75
$ws.Range("N46").Value = -2877.6155
$ws.Range("H122").Value = 100004170
$ws.Range("I122").Value = 142861330
$ws.Range("K122").Value = 428583990
$ws.Range("M122").Value = -428581540
$ws.Range("H126").Value = 3069.64
$ws.Range("I126").Value = 2176.9285
$ws.Range("J126").Value = 4205.8184
$ws.Range("K126").Value = 6530.7855
$ws.Range("L126").Value = 12617.4552
$ws.Range("M126").Value = -4060.7855
$ws.Range("N126").Value = -17557.4552
$ws.Range("H132").Value = 1836
$ws.Range("I132").Value = 1780.6666
$ws.Range("K132").Value = 5341.9998
$ws.Range("M132").Value = -2811.9998
$ws.Range("H136").Value = 3491.5
$ws.Range("I136").Value = 4457.364
$ws.Range("K136").Value = 13372.092
$ws.Range("M136").Value = -10822.092

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 5000
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H41").Value = 9830.6
$ws.Range("J41").Value = 9733
$ws.Range("L41").Value = 9733
$ws.Range("N41").Value = -10513
$ws.Range("H45").Value = 25518
$ws.Range("I45").Value = 15000
$ws.Range("J45").Value = 30777
$ws.Range("K45").Value = 15000
$ws.Range("L45").Value = 30777
$ws.Range("M45").Value = -14509
$ws.Range("N45").Value = -31759
$ws.Range("H126").Value = 2905.6155
$ws.Range("I126").Value = 2370.353
$ws.Range("K126").Value = 7111.059
$ws.Range("M126").Value = -4641.059
$ws.Range("H132").Value = 2149.9375
$ws.Range("I132").Value = 1761.4615
$ws.Range("J132").Value = 3833.3333
$ws.Range("K132").Value = 5284.3845
$ws.Range("L132").Value = 11499.9999
$ws.Range("M132").Value = -2754.3845
$ws.Range("N132").Value = -16559.9999
$ws.Range("H136").Value = 1981.5
$ws.Range("I136").Value = 1773.4375
$ws.Range("J136").Value = 2536.3333
$ws.Range("K136").Value = 5320.3125
$ws.Range("L136").Value = 7608.999899999999
$ws.Range("M136").Value = -2770.3125
$ws.Range("N136").Value = -12708.9999
